# Daily attendance processing - swap "System, <email>" to "<email>, System"
# in the "Recorded By" column (column G) of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value()

    if ($val -ne $null -and $val -like "System, *") {
        $parts = $val -split ", "
        if ($parts.Count -eq 2 -and $parts[0] -eq "System") {
            $cell.Value = "$($parts[1]), System"
        }
    }
}
